$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 14 (pushes old rows 14-18 down to 15-19),
#    making room for the new item "ريكسونه رجالى" as item #8.
$ws.Rows("14:14").Insert()

# 2. Clone the formatting of the row directly below (old row 14, now row 15,
#    which still carries the original item-row style pattern) into the new
#    blank row 14.
$ws.Range("A15:Q15").Copy($ws.Range("A14:Q14"))
$ws.Rows("14:14").RowHeight = 25.5

# 3. Re-create the merged cells for the new row 14 (Insert() does not clone
#    merges into the freshly inserted row).
$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# 4. Fill in the new item's data: رقم 8, ريكسونه رجالى, رصيد 4:0, سعر 27.00, سعر بيع 27.0000
$ws.Range("A14").Value2 = 8
$ws.Range("C14").Value2 = "ريكسونه رجالى"
$ws.Range("H14").Value2 = "4:0"
$ws.Range("L14").Value2 = "0"
$ws.Range("N14").Value2 = "27.00"
$ws.Range("P14").Value2 = "27.0000"
$ws.Range("Q14").Value2 = "1:0"

# 5. Update the totals row (old row 17, now shifted to row 18): add the new
#    item's price (27.00) to the previous grand total.
$ws.Range("P18").Value2 = 369.72000000000003
$ws.Rows("18:18").RowHeight = 24.75

# 6. Update the generated timestamp footer (old row 18, now row 19) to match
#    the new export time.
$ws.Range("A19").Value2 = "Monday, 18 August, 2025 9:57 AM"

Write-Host "Edit complete"
